$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 827.8182
$ws.Range("I32").Value = 312.5
$ws.Range("J32").Value = 1122.2858
$ws.Range("K32").Value = 312.5
$ws.Range("L32").Value = 1122.2858
$ws.Range("M32").Value = 13.5
$ws.Range("N32").Value = -1774.2858

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 1925
$ws.Range("I51").Value = 900
$ws.Range("J51").Value = 2266.6667
$ws.Range("K51").Value = 900
$ws.Range("L51").Value = 2266.6667
$ws.Range("M51").Value = -416
$ws.Range("N51").Value = -3234.6667

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H88").Value = 823480.7
$ws.Range("J88").Value = 3086793.2
$ws.Range("L88").Value = 3086793.2
$ws.Range("N88").Value = -3087605.2

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H91").Value = 823480.7
$ws.Range("J91").Value = 3086793.2
$ws.Range("L91").Value = 3086793.2
$ws.Range("N91").Value = -3089601.2

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 7572.0527
$ws.Range("I106").Value = 7826.0557
$ws.Range("K106").Value = 7826.0557
$ws.Range("M106").Value = -7195.0557

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1369.8
$ws.Range("I137").Value = 1220.2439
$ws.Range("J137").Value = 1625.2916
$ws.Range("K137").Value = 3660.7317
$ws.Range("L137").Value = 4875.8748
$ws.Range("M137").Value = -1110.7317
$ws.Range("N137").Value = -9975.8748

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2031.6768
$ws.Range("I138").Value = 1361.3928
$ws.Range("J138").Value = 2296.0142
$ws.Range("K138").Value = 4084.1784
$ws.Range("L138").Value = 6888.042600000001
$ws.Range("M138").Value = 1055.8216
$ws.Range("N138").Value = -17168.0426

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H133").Value = 34220
$ws.Range("J133").Value = 34220
$ws.Range("L133").Value = 34220
$ws.Range("N133").Value = -39280

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2627.9
$ws.Range("I20").Value = 2698.7778
$ws.Range("J20").Value = 1990
$ws.Range("K20").Value = 2698.7778
$ws.Range("L20").Value = 1990
$ws.Range("M20").Value = -2451.7778
$ws.Range("N20").Value = -2484

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 58824830
$ws.Range("I105").Value = 66667744
$ws.Range("K105").Value = 66667744
$ws.Range("M105").Value = -66665997

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 76924584
$ws.Range("I16").Value = 100001540
$ws.Range("J16").Value = 1416.6666
$ws.Range("K16").Value = 100001540
$ws.Range("L16").Value = 1416.6666
$ws.Range("M16").Value = -100001253
$ws.Range("N16").Value = -1990.6666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1384.0754
$ws.Range("I31").Value = 1374.9756
$ws.Range("J31").Value = 1415.1666
$ws.Range("K31").Value = 1374.9756
$ws.Range("L31").Value = 1415.1666
$ws.Range("M31").Value = -1079.9756
$ws.Range("N31").Value = -2005.1666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 1384.0754
$ws.Range("I34").Value = 1374.9756
$ws.Range("J34").Value = 1415.1666
$ws.Range("K34").Value = 1374.9756
$ws.Range("L34").Value = 1415.1666
$ws.Range("M34").Value = -1172.9756
$ws.Range("N34").Value = -1819.1666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 8105
$ws.Range("I58").Value = 1483.2
$ws.Range("J58").Value = 10651.846
$ws.Range("K58").Value = 1483.2
$ws.Range("L58").Value = 10651.846
$ws.Range("M58").Value = -1280.2
$ws.Range("N58").Value = -11057.846

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 1630.6875
$ws.Range("I99").Value = 1564.6364
$ws.Range("J99").Value = 1776
$ws.Range("K99").Value = 1564.6364
$ws.Range("L99").Value = 1776
$ws.Range("M99").Value = -66.63640000000009
$ws.Range("N99").Value = -4772

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H112").Value = 34221.89
$ws.Range("J112").Value = 34221.89
$ws.Range("L112").Value = 34221.89
$ws.Range("N112").Value = -37175.89

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 76924584
$ws.Range("I113").Value = 100001540
$ws.Range("J113").Value = 1416.6666
$ws.Range("K113").Value = 100001540
$ws.Range("L113").Value = 1416.6666
$ws.Range("M113").Value = -99999370
$ws.Range("N113").Value = -5756.6666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 1630.6875
$ws.Range("I126").Value = 1564.6364
$ws.Range("J126").Value = 1776
$ws.Range("K126").Value = 4693.9092
$ws.Range("L126").Value = 5328
$ws.Range("M126").Value = -2223.9092
$ws.Range("N126").Value = -10268

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2608.3572
$ws.Range("I132").Value = 1947.8
$ws.Range("J132").Value = 2975.3333
$ws.Range("K132").Value = 5843.4
$ws.Range("L132").Value = 8925.999899999999
$ws.Range("M132").Value = -3313.4
$ws.Range("N132").Value = -13985.9999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 8105
$ws.Range("I136").Value = 1483.2
$ws.Range("J136").Value = 10651.846
$ws.Range("K136").Value = 4449.6
$ws.Range("L136").Value = 31955.538
$ws.Range("M136").Value = -1899.6
$ws.Range("N136").Value = -37055.538

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H63").Value = 15058.9
$ws.Range("I63").Value = 3000
$ws.Range("J63").Value = 20227
$ws.Range("K63").Value = 9000
$ws.Range("L63").Value = 60681
$ws.Range("M63").Value = -8251
$ws.Range("N63").Value = -62179

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H64").Value = 2947.1
$ws.Range("I64").Value = 1064.6666
$ws.Range("J64").Value = 3753.8572
$ws.Range("K64").Value = 3193.9998
$ws.Range("L64").Value = 11261.5716
$ws.Range("M64").Value = -2923.9998
$ws.Range("N64").Value = -11801.5716

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H66").Value = 15058.9
$ws.Range("I66").Value = 3000
$ws.Range("J66").Value = 20227
$ws.Range("K66").Value = 27000
$ws.Range("L66").Value = 182043
$ws.Range("M66").Value = -23256
$ws.Range("N66").Value = -189531

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H67").Value = 2947.1
$ws.Range("I67").Value = 1064.6666
$ws.Range("J67").Value = 3753.8572
$ws.Range("K67").Value = 3193.9998
$ws.Range("L67").Value = 11261.5716
$ws.Range("M67").Value = -2257.9998
$ws.Range("N67").Value = -13133.5716

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 628.34784
$ws.Range("I113").Value = 529.05884
$ws.Range("J113").Value = 686.5517
$ws.Range("K113").Value = 1587.17652
$ws.Range("L113").Value = 2059.6551
$ws.Range("M113").Value = 582.82348
$ws.Range("N113").Value = -6399.6551

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H133").Value = 46422.5
$ws.Range("J133").Value = 46422.5
$ws.Range("L133").Value = 46422.5
$ws.Range("N133").Value = -51482.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 2244.889
$ws.Range("I136").Value = 2029.1428
$ws.Range("K136").Value = 6087.428400000001
$ws.Range("M136").Value = -3537.428400000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2069.9
$ws.Range("I132").Value = 1948.8182
$ws.Range("J132").Value = 2640.7144
$ws.Range("K132").Value = 5846.4546
$ws.Range("L132").Value = 7922.1432
$ws.Range("M132").Value = -3316.4546
$ws.Range("N132").Value = -12982.1432
